$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) values stay text, matching the source format
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '26.779.49'
$ws.Range("E2").Value = '  +0.66%  '
$ws.Range("D3").Value = '1.732.31'
$ws.Range("E3").Value = '  -0.40%  '
$ws.Range("D4").Value = '0.9961'
$ws.Range("E4").Value = '  -0.42%  '
$ws.Range("D5").Value = '241.97'
$ws.Range("E5").Value = '  -1.62%  '
$ws.Range("D6").Value = '0.9966'
$ws.Range("E6").Value = '  -0.39%  '
$ws.Range("D7").Value = '0.4945'
$ws.Range("E7").Value = '  +0.61%  '
$ws.Range("D8").Value = '0.2617'
$ws.Range("E8").Value = '  -2.11%  '
$ws.Range("D9").Value = '0.06224'
$ws.Range("E9").Value = '  -0.79%  '
$ws.Range("D10").Value = '1.728.27'
$ws.Range("E10").Value = '  -0.74%  '
$ws.Range("D11").Value = '15.76'
$ws.Range("E11").Value = '  +0.24%  '
$ws.Range("D12").Value = '0.06980'
$ws.Range("E12").Value = '  -0.90%  '
$ws.Range("D13").Value = '0.6135'
$ws.Range("E13").Value = '  +0.03%  '
$ws.Range("D14").Value = '4.499'
$ws.Range("E14").Value = '  -1.79%  '
$ws.Range("D15").Value = '77.16'
$ws.Range("E15").Value = '  -1.05%  '
$ws.Range("D16").Value = '0.9962'
$ws.Range("E16").Value = '  -0.42%  '
$ws.Range("D17").Value = '26.543.74'
$ws.Range("E17").Value = '  -0.28%  '
$ws.Range("D18").Value = '0.9963'
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("D19").Value = '0.000007172'
$ws.Range("E19").Value = '  -1.13%  '
$ws.Range("D20").Value = '11.41'
$ws.Range("E20").Value = '  -1.33%  '
$ws.Range("D21").Value = '1.947.79'
$ws.Range("E21").Value = '  -1.29%  '
$ws.Range("D22").Value = '4.439'
$ws.Range("E22").Value = '  -2.65%  '
$ws.Range("D23").Value = '8.517'
$ws.Range("E23").Value = '  -2.22%  '
$ws.Range("D24").Value = '5.123'
$ws.Range("E24").Value = '  -2.83%  '
$ws.Range("D25").Value = '139.00'
$ws.Range("E25").Value = '  -0.06%  '
$ws.Range("D26").Value = '15.34'
$ws.Range("E26").Value = '  -0.48%  '
$ws.Range("D27").Value = '1.413'
$ws.Range("E27").Value = '  -0.37%  '
$ws.Range("D28").Value = '1.760'
$ws.Range("E28").Value = '  -0.07%  '
$ws.Range("D29").Value = '106.41'
$ws.Range("E29").Value = '  -0.90%  '
$ws.Range("D30").Value = '3.940'
$ws.Range("E30").Value = '  -1.98%  '
$ws.Range("D31").Value = '0.07979'
$ws.Range("D32").Value = '3.659'
$ws.Range("E32").Value = '  -1.72%  '
$ws.Range("D33").Value = '0.04523'
$ws.Range("E33").Value = '  -1.95%  '
$ws.Range("B34").Value = 'HuobiToken'
$ws.Range("C34").Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range("D34").Value = '2.608'
$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("B35").Value = 'ARBITRUM'
$ws.Range("C35").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D35").Value = '1.003'
$ws.Range("E35").Value = '  -0.99%  '
$ws.Range("B36").Value = 'ImmutableX'
$ws.Range("C36").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D36").Value = '0.6272'
$ws.Range("E36").Value = '  -1.88%  '
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").Value = '0.9452'
$ws.Range("E37").Value = '  +4.45%  '
$ws.Range("B38").Value = 'RenderToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D38").Value = '2.016'
$ws.Range("E38").Value = '  -2.20%  '
$ws.Range("B39").Value = 'MXToken'
$ws.Range("C39").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D39").Value = '2.410'
$ws.Range("E39").Value = '  -0.71%  '
$ws.Range("B40").Value = 'PaxDollar'
$ws.Range("C40").Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range("D40").Value = '0.9964'
$ws.Range("E40").Value = '  -0.65%  '
$ws.Range("B41").Value = 'VeChain'
$ws.Range("C41").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D41").Value = '0.01504'
$ws.Range("E41").Value = '  +0.04%  '
$ws.Range("B42").Value = 'Quant'
$ws.Range("C42").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D42").Value = '99.87'
$ws.Range("E42").Value = '  -2.19%  '
$ws.Range("B43").Value = 'FraxShare'
$ws.Range("C43").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D43").Value = '5.480'
$ws.Range("E43").Value = '  +0.91%  '
$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '0.3866'
$ws.Range("E44").Value = '  -1.52%  '
$ws.Range("B45").Value = 'Aptos'
$ws.Range("C45").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$ws.Range("D45").Value = '6.945'
$ws.Range("E45").Value = '  +1.13%  '
$ws.Range("B46").Value = 'Algorand'
$ws.Range("C46").Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range("D46").Value = '0.1163'
$ws.Range("E46").Value = '  -1.85%  '
$ws.Range("B47").Value = 'Cronos'
$ws.Range("C47").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D47").Value = '0.05395'
$ws.Range("E47").Value = '  +0.00%  '
$ws.Range("B48").Value = 'Elrond'
$ws.Range("C48").Value = 'https://coinranking.com/coin/omwkOTglq+elrond-egld'
$ws.Range("D48").Value = '30.52'
$ws.Range("E48").Value = '  -0.18%  '
$ws.Range("B49").Value = 'EnergySwap'
$ws.Range("C49").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D49").Value = '7.765'
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").Value = '51.74'
$ws.Range("E50").Value = '  -0.05%  '
$ws.Range("B51").Value = 'NEARProtocol'
$ws.Range("C51").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D51").Value = '1.227'
$ws.Range("E51").Value = '  -2.26%  '
